$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" date column (C) for rows 2-9 from 45224 to 45233
foreach ($row in 2..9) {
    $ws.Cells.Item($row, 3).Value = 45233
}
